# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.2881169905109251, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 0, 1.84748871573303)
    3 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248)
    4 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 3.755628166162433)
    5 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
